$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 188.22223
$ws.Range("I41").Value = 188.22223
$ws.Range("K41").Value = 188.22223
$ws.Range("M41").Value = 251.77777

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 5799.2
$ws.Range("I43").Value = 8332
$ws.Range("J43").Value = 2000
$ws.Range("K43").Value = 8332
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = -8263
$ws.Range("N43").Value = -2138

# Row 47 (Leve Item ID 2169)
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

# Row 61 (Leve Item ID 4604)
$ws.Range("H61").Value = 3008
$ws.Range("I61").Value = 999
$ws.Range("J61").Value = 5017
$ws.Range("K61").Value = 2997
$ws.Range("L61").Value = 15051
$ws.Range("M61").Value = -2825
$ws.Range("N61").Value = -15395

# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 2088.889
$ws.Range("I64").Value = 1200
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 1200
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -952
$ws.Range("N64").Value = -3696

# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 2088.889
$ws.Range("I67").Value = 1200
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 1200
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -342
$ws.Range("N67").Value = -4916

# Row 81 (Leve Item ID 10637)
$ws.Range("H81").Value = 63999
$ws.Range("J81").Value = 63999
$ws.Range("L81").Value = 63999
$ws.Range("N81").Value = -65995

# Row 84 (Leve Item ID 10637)
$ws.Range("H84").Value = 63999
$ws.Range("J84").Value = 63999
$ws.Range("L84").Value = 191997
$ws.Range("N84").Value = -201981

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 2166.1667
$ws.Range("I86").Value = 498.5
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 498.5
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = 624.5
$ws.Range("N86").Value = -5246

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 2166.1667
$ws.Range("I89").Value = 498.5
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 2492.5
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = 3123.5
$ws.Range("N89").Value = -26232

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1389.1052
$ws.Range("I98").Value = 1456.1333
$ws.Range("J98").Value = 1137.75
$ws.Range("K98").Value = 1456.1333
$ws.Range("L98").Value = 1137.75
$ws.Range("M98").Value = 41.86670000000004
$ws.Range("N98").Value = -4133.75

# Row 101 (Leve Item ID 19884)
$ws.Range("H101").Value = 641.5
$ws.Range("I101").Value = 641.5
$ws.Range("K101").Value = 1924.5
$ws.Range("M101").Value = -302.5

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 2000
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -8884

# Row 117 (Leve Item ID 26118)
$ws.Range("H117").Value = 100000
$ws.Range("J117").Value = 100000
$ws.Range("L117").Value = 100000
$ws.Range("N117").Value = -109178

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1389.1052
$ws.Range("I122").Value = 1456.1333
$ws.Range("J122").Value = 1137.75
$ws.Range("K122").Value = 4368.3999
$ws.Range("L122").Value = 3413.25
$ws.Range("M122").Value = -1918.3999
$ws.Range("N122").Value = -8313.25

# Row 131 (Leve Item ID 36108)
$ws.Range("H131").Value = 813.8570999999999
$ws.Range("I131").Value = 682.8333
$ws.Range("J131").Value = 1600
$ws.Range("K131").Value = 2048.4999
$ws.Range("L131").Value = 4800
$ws.Range("M131").Value = 2991.5001
$ws.Range("N131").Value = -14880

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2315.6667
$ws.Range("I132").Value = 2473.75
$ws.Range("K132").Value = 7421.25
$ws.Range("M132").Value = -4891.25

$ws = $wb.Worksheets.Item("ARM")
# Row 31 (Leve Item ID 19533)
$ws.Range("H31").Value = 2885.5
$ws.Range("I31").Value = 2885.5
$ws.Range("K31").Value = 2885.5
$ws.Range("M31").Value = -2591.5

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 1061.4
$ws.Range("I88").Value = 433.33334
$ws.Range("K88").Value = 433.33334
$ws.Range("M88").Value = -27.33334000000002

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 1061.4
$ws.Range("I91").Value = 433.33334
$ws.Range("K91").Value = 433.33334
$ws.Range("M91").Value = 970.66666

# Row 95 (Leve Item ID 18204)
$ws.Range("H95").Value = 48333
$ws.Range("J95").Value = 48333
$ws.Range("L95").Value = 48333
$ws.Range("N95").Value = -53825

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1392.4286
$ws.Range("I132").Value = 1407.1666
$ws.Range("J132").Value = 1304
$ws.Range("K132").Value = 4221.4998
$ws.Range("L132").Value = 3912
$ws.Range("M132").Value = -1691.4998
$ws.Range("N132").Value = -8972

$ws = $wb.Worksheets.Item("BSM")
# Row 82 (Leve Item ID 11877)
$ws.Range("H82").Value = 154017.38
$ws.Range("I82").Value = 6450.8
$ws.Range("J82").Value = 399961.66
$ws.Range("K82").Value = 6450.8
$ws.Range("L82").Value = 399961.66
$ws.Range("M82").Value = -6067.8
$ws.Range("N82").Value = -400727.66

# Row 85 (Leve Item ID 11877)
$ws.Range("H85").Value = 154017.38
$ws.Range("I85").Value = 6450.8
$ws.Range("J85").Value = 399961.66
$ws.Range("K85").Value = 6450.8
$ws.Range("L85").Value = 399961.66
$ws.Range("M85").Value = -5124.8
$ws.Range("N85").Value = -402613.66

# Row 102 (Leve Item ID 19565)
$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 10000
$ws.Range("K102").Value = 10000
$ws.Range("M102").Value = -6755

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2358.0952
$ws.Range("I134").Value = 1551.4286
$ws.Range("K134").Value = 4654.2858
$ws.Range("M134").Value = -2119.2858

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 89.111115
$ws.Range("I7").Value = 38.714287
$ws.Range("J7").Value = 265.5
$ws.Range("K7").Value = 38.714287
$ws.Range("L7").Value = 265.5
$ws.Range("M7").Value = 74.285713
$ws.Range("N7").Value = -491.5

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3803.4211
$ws.Range("I31").Value = 2756.4614
$ws.Range("J31").Value = 6071.8335
$ws.Range("K31").Value = 2756.4614
$ws.Range("L31").Value = 6071.8335
$ws.Range("M31").Value = -2461.4614
$ws.Range("N31").Value = -6661.8335

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3803.4211
$ws.Range("I34").Value = 2756.4614
$ws.Range("J34").Value = 6071.8335
$ws.Range("K34").Value = 2756.4614
$ws.Range("L34").Value = 6071.8335
$ws.Range("M34").Value = -2554.4614
$ws.Range("N34").Value = -6475.8335

# Row 50 (Leve Item ID 1862)
$ws.Range("H50").Value = 28379.6
$ws.Range("J50").Value = 28379.6
$ws.Range("L50").Value = 28379.6
$ws.Range("N50").Value = -29629.6

# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 82126.60000000001
$ws.Range("I62").Value = 2658.5
$ws.Range("K62").Value = 2658.5
$ws.Range("M62").Value = -2034.5

# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 82126.60000000001
$ws.Range("I65").Value = 2658.5
$ws.Range("K65").Value = 13292.5
$ws.Range("M65").Value = -10172.5

# Row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 57997
$ws.Range("J74").Value = 57997
$ws.Range("L74").Value = 57997
$ws.Range("N74").Value = -59745

# Row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 57997
$ws.Range("J77").Value = 57997
$ws.Range("L77").Value = 173991
$ws.Range("N77").Value = -182727

# Row 97 (Leve Item ID 19730)
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 316.2
$ws.Range("I2").Value = 8
$ws.Range("K2").Value = 8
$ws.Range("M2").Value = 105

# Row 38 (Leve Item ID 2737)
$ws.Range("H38").Value = 15000
$ws.Range("J38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("N38").Value = -15926

# Row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 17999.8
$ws.Range("J46").Value = 20000
$ws.Range("L46").Value = 20000
$ws.Range("N46").Value = -20312

# Row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 44965
$ws.Range("J57").Value = 44956.25
$ws.Range("L57").Value = 44956.25
$ws.Range("N57").Value = -46596.25

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 36330.69
$ws.Range("I122").Value = 1317.5834
$ws.Range("J122").Value = 204393.6
$ws.Range("K122").Value = 3952.7502
$ws.Range("L122").Value = 613180.8
$ws.Range("M122").Value = -1502.7502
$ws.Range("N122").Value = -618080.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2631.7273
$ws.Range("I7").Value = 2438.7778
$ws.Range("K7").Value = 2438.7778
$ws.Range("M7").Value = -2326.7778

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 4162
$ws.Range("I22").Value = 2256.25
$ws.Range("J22").Value = 7973.5
$ws.Range("K22").Value = 2256.25
$ws.Range("L22").Value = 7973.5
$ws.Range("M22").Value = -1961.25
$ws.Range("N22").Value = -8563.5

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 4162
$ws.Range("I27").Value = 2256.25
$ws.Range("J27").Value = 7973.5
$ws.Range("K27").Value = 2256.25
$ws.Range("L27").Value = 7973.5
$ws.Range("M27").Value = -2149.25
$ws.Range("N27").Value = -8187.5

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 904.1429000000001
$ws.Range("I55").Value = 590
$ws.Range("K55").Value = 590
$ws.Range("M55").Value = -417

# Row 81 (Leve Item ID 10897)
$ws.Range("H81").Value = 4482
$ws.Range("J81").Value = 7800
$ws.Range("L81").Value = 7800
$ws.Range("N81").Value = -9796

# Row 84 (Leve Item ID 10897)
$ws.Range("H84").Value = 4482
$ws.Range("J84").Value = 7800
$ws.Range("L84").Value = 23400
$ws.Range("N84").Value = -33384

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2631.7273
$ws.Range("I126").Value = 2438.7778
$ws.Range("K126").Value = 7316.3334
$ws.Range("M126").Value = -4846.3334

# Row 131 (Leve Item ID 35466)
$ws.Range("H131").Value = 22000
$ws.Range("J131").Value = 22000
$ws.Range("L131").Value = 22000
$ws.Range("N131").Value = -32080

$ws = $wb.Worksheets.Item("WVR")
# Row 56 (Leve Item ID 10912)
$ws.Range("H56").Value = 28649.75
$ws.Range("I56").Value = 21428.334
$ws.Range("K56").Value = 21428.334
$ws.Range("M56").Value = -20714.334

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 1775
$ws.Range("J107").Value = 1699
$ws.Range("L107").Value = 5097
$ws.Range("N107").Value = -8937

